$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updates to existing rows (row 2, 3, 6, 7) ---

# Row 2
$ws.Range("K2").Value = 13

# Row 3
$ws.Range("J3").Value = 1.01
$ws.Range("L3").Value = 1.1

# Row 6
$ws.Range("H6").Value = 3.4
$ws.Range("I6").Value = 2.15
$ws.Range("U6").Value = 12
$ws.Range("V6").Value = 9
$ws.Range("W6").Value = 26
$ws.Range("X6").Value = 19
$ws.Range("AA6").Value = 5.8
$ws.Range("AB6").Value = 11.75
$ws.Range("AE6").Value = 6.9
$ws.Range("AF6").Value = 9
$ws.Range("AH6").Value = 16.5
$ws.Range("AJ6").Value = 22

# Row 7
$ws.Range("G7").Value = 2.35
$ws.Range("H7").Value = 3.3
$ws.Range("I7").Value = 2.62
$ws.Range("P7").Value = 1.37
$ws.Range("Q7").Value = 2.5
$ws.Range("T7").Value = 6.6
$ws.Range("U7").Value = 9.5
$ws.Range("W7").Value = 19
$ws.Range("X7").Value = 16
$ws.Range("Z7").Value = 9.25
$ws.Range("AA7").Value = 5.6
$ws.Range("AD7").Value = 400
$ws.Range("AE7").Value = 7
$ws.Range("AF7").Value = 10.75
$ws.Range("AH7").Value = 23

# --- New rows 8-11 ---

$newRows = @(
    @("YaJXV68c", "17/06/2025", "07:30", "SOUTH KOREA - K LEAGUE 1", "Daegu", "Pohang", 3.8, 3.5, 1.95, 1.05, 11, 1.25, 3.75, 1.83, 1.98, 1.36, 3, 1.67, 2.1, 12, 21, 13, 41, 29, 34, 11, 6.5, 13, 41, 201, 8, 9.5, 8.5, 17, 15, 23),
    @("KKyo83Gj", "17/06/2025", "07:30", "SOUTH KOREA - K LEAGUE 1", "Jeonbuk", "Suwon FC", 1.6, 3.9, 5.5, 1.03, 13, 1.22, 3.75, 1.8, 2, 1.33, 3.25, 1.8, 1.8, 7.5, 7.5, 8.5, 12, 13, 26, 11, 7.5, 17, 51, 251, 15, 29, 17, 51, 41, 41),
    @("EL9DQBWH", "17/06/2025", "07:30", "SOUTH KOREA - K LEAGUE 1", "Seoul", "Gangwon", 1.67, 3.6, 5.25, 1.07, 9, 1.36, 3, 2.2, 1.65, 1.44, 2.63, 2.1, 1.63, 5.5, 7, 8.5, 12, 15, 34, 8, 7, 21, 67, 501, 11, 26, 17, 51, 41, 51),
    @("IkzEmucm", "17/06/2025", "20:00", "URUGUAY - LIGA AUF URUGUAYA", "Wanderers", "Defensor Sp.", 2.6, 3.1, 2.7, 1.07, 9, 1.4, 2.75, 2.25, 1.62, 1.44, 2.63, 1.91, 1.8, 7.5, 12, 11, 26, 23, 34, 8, 6, 15, 51, 351, 8, 13, 11, 29, 23, 34)
)

$rowIndex = 8
foreach ($rowData in $newRows) {
    $colIndex = 1
    foreach ($value in $rowData) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $value
        $colIndex++
    }
    $rowIndex++
}
